# Fix Training Data Issue (#48)
# The "Date" column (BF) held the wrong literal date string - it was off
# by one day because of the way NBA stats dates were originally captured
# ("6-14-2012-13" instead of the proper ISO date "2013-06-14"). Update the
# 30 data rows (BF2:BF31) to the corrected date text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-14-2012-13"
$newDate = "2013-06-14"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value() -eq $oldDate) {
        # Force text interpretation so Excel doesn't auto-convert the
        # "2013-06-14" literal into a date serial value, then clear the
        # temporary number-format override so the cell keeps its original
        # (default) style.
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
        $cell.ClearFormats()
    }
}
